# Generate Report for Handoff
# A new handoff run updates the "Latest Handoff Date" / "Latest Handoff
# Datetime" timestamp for every file row that is mid-flight (Handback
# transform failed / Ready for handoff) on each of the report's sheets.
# Rows that are already synced ("Handed back: in sync with en-US") or are
# still "In Translation" keep their prior timestamp.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D ("Latest Handoff Date")
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $overviewRows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-30 11:14:40"
}

# zh-cn sheet: column E ("Latest Handoff Datetime")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $zhCnRows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-30 11:14:28"
}

# de-de sheet: column E ("Latest Handoff Datetime")
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $deDeRows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-30 11:14:40"
}
